$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: append a new sentence to the paragraph that ends with
# "...design and optimization of digital products. " :
#   "I have also completed projects involving data analysis, machine
#   learning and communication of insights."
# split across 8 separate runs (matching the target OOXML), each
# carrying the same Times New Roman / 24 / en-US run formatting as
# the rest of the paragraph.
# -----------------------------------------------------------------

$find = $d.Content
$find.Find.Execute("digital products.") | Out-Null
$para = $find.Paragraphs(1)
$pEnd = $para.Range.End
$insPos = $pEnd - 1   # just before the paragraph mark

$segments = @(
    "I have ",
    "also ",
    "completed projects ",
    "involving",
    " data analysis, machine learning and ",
    "communication",
    " of insights",
    "."
)

$full = ""
foreach ($seg in $segments) { $full = $full + $seg }

$startPos = $insPos

# Insert the whole block at once via a collapsed range right before the
# paragraph mark -- this makes the new text inherit the run formatting
# (font / size / language) of the text immediately to its left.
$insertRange = $d.Range($insPos, $insPos)
$insertRange.Text = $full
$endPos = $startPos + $full.Length

# Compute the offset of every segment boundary (including the very
# first one, which separates the new text from the pre-existing
# trailing-space run).
$offsets = @($startPos)
$cur = $startPos
foreach ($seg in $segments) {
    $cur = $cur + $seg.Length
    $offsets += $cur
}

# Force a run break at each boundary (except the very last, which is
# just the end of the inserted text) by toggling Bold on/off -- this
# is a no-op visually but makes the engine split the run there while
# keeping/re-deriving the correct rPr on both sides.
for ($i = 0; $i -lt $offsets.Length - 1; $i++) {
    $b = $offsets[$i]
    $splitRange = $d.Range($b, $endPos)
    $splitRange.Bold = 1
    $splitRange.Bold = 0
}

# -----------------------------------------------------------------
# Change 2: split "I also have experience with event management and
# facilitating general assemblies. " into two runs:
#   "Furthermore, I " + "have experience with event management and
#   facilitating general assemblies. "
# -----------------------------------------------------------------

$find2 = $d.Content
$find2.Find.Execute("I also have experience with event management") | Out-Null
Write-Host "Found second target: $($find2.Find.Found)"

$r2 = $d.Content
$r2.Find.Execute("I also ") | Out-Null
$r2.Text = "Furthermore, I "
$r2.Bold = 1
$r2.Bold = 0

Write-Host "Done applying edits"
